$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Cells.Item(20,1)
try {
  $s = $r.Style
  Write-Output "Style: $s"
} catch {
  Write-Output "ERROR: $_"
}
